$wb = $excel.ActiveWorkbook
$odiBatting = $wb.Worksheets.Item("ODI Batting")

# --- 1) Add new worksheet "ODI Batting Extra" after "ODI Bowling" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Match the page margins used by the other sheets in the workbook
# (0.75in/0.75in/1in/1in/0.5in/0.5in -- Excel's PageSetup uses points).
$ps = $newSheet.PageSetup
$ps.LeftMargin = 54
$ps.RightMargin = 54
$ps.TopMargin = 72
$ps.BottomMargin = 72
$ps.HeaderMargin = 36
$ps.FooterMargin = 36

# Copy header formatting (bold/border/center) from an existing header row
$odiBatting.Range("A1:F1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122) # xlPasteFormats

# Header row values
$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# MATCH_CODE values should be text (not numbers). Copy the *values* from
# existing text cells elsewhere in the workbook so the new cells come out
# text-typed without registering any new number-format/style.
$odiBatting.Range("D2").Copy()
$newSheet.Range("A2").PasteSpecial(-4163) # xlPasteValues
$odiBatting.Range("D3").Copy()
$newSheet.Range("A3").PasteSpecial(-4163)
$odiBatting.Range("D4").Copy()
$newSheet.Range("A4").PasteSpecial(-4163)

# BATTING_POSITION is numeric
$newSheet.Range("B2").Value = 10
$newSheet.Range("B3").Value = 10
$newSheet.Range("B4").Value = 10

# NUM_4 / NUM_6 for the 3rd data row are text "0" -- reuse an existing
# text "0" cell's value the same way.
$odiBatting.Range("I4").Copy()
$newSheet.Range("C4").PasteSpecial(-4163)
$odiBatting.Range("I4").Copy()
$newSheet.Range("D4").PasteSpecial(-4163)

# MAN_OF_MATCH column
$newSheet.Range("F2").Value = "NO"
$newSheet.Range("F3").Value = "NO"
$newSheet.Range("F4").Value = "NO"

$excel.CutCopyMode = $false

# --- 2) Remove the now-empty B2/B3 cells on "ODI Batting" sheet ---
$odiBatting.Range("B2").ClearContents()
$odiBatting.Range("B3").ClearContents()

# Restore the workbook's original active sheet (adding a sheet makes it active)
$wb.Worksheets.Item("Player Info").Activate()
